$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push the existing "Et1"/"dos" data rows down by two, opening up a blank
# row 2 and a row 3 that will hold the new "ORG" line.
$ws.Rows.Item(2).Insert()
$ws.Rows.Item(2).Insert()

# New ORG row (row 3) - column A stays blank.
$ws.Cells.Item(3, 2).Value = "ORG"
$ws.Cells.Item(3, 3).Value = "%00001111"

# Three new rows appended after the existing data (rows 6-8).
$ws.Cells.Item(6, 2).Value = "SWI"

$ws.Cells.Item(7, 2).Value = "DS.B"
$ws.Cells.Item(7, 3).Value = "%0011000011111100"

$ws.Cells.Item(8, 2).Value = "END"
